$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values are treated as text (many look numeric, e.g. "1.003")
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.454.39'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").Value = '1.890.39'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '238.64'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = '0.4842'
$ws.Range("E7").Value = '  -1.42%  '
$ws.Range("D8").Value = '0.2896'
$ws.Range("E8").Value = '  -1.55%  '
$ws.Range("D9").Value = '0.06606'
$ws.Range("E9").Value = '  -1.47%  '
$ws.Range("D10").Value = '1.890.23'
$ws.Range("E10").Value = '  -0.17%  '
$ws.Range("D11").Value = '16.80'
$ws.Range("E11").Value = '  -1.11%  '
$ws.Range("D12").Value = '0.07407'
$ws.Range("E12").Value = '  +0.90%  '
$ws.Range("D13").Value = '5.185'
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("D14").Value = '88.60'
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").Value = '0.6616'
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").Value = '30.430.61'
$ws.Range("E16").Value = '  -0.64%  '
$ws.Range("D17").Value = '13.53'
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("D18").Value = '0.000007770'
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").Value = '2.143.82'
$ws.Range("E20").Value = '  +0.22%  '
$ws.Range("D21").Value = '5.367'
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '221.34'
$ws.Range("E23").Value = '  +16.02%  '
$ws.Range("D24").Value = '6.220'
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").Value = '9.376'
$ws.Range("E25").Value = '  -1.69%  '
$ws.Range("D26").Value = '163.05'
$ws.Range("E26").Value = '  +0.84%  '
$ws.Range("D27").Value = '18.86'
$ws.Range("E27").Value = '  +2.04%  '
$ws.Range("D28").Value = '1.942'
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("D29").Value = '1.450'
$ws.Range("E29").Value = '  -1.08%  '
$ws.Range("D30").Value = '4.331'
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("D31").Value = '0.09215'
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("D32").Value = '4.038'
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").Value = '0.05074'
$ws.Range("E33").Value = '  -3.23%  '
$ws.Range("D34").Value = '0.7595'
$ws.Range("D35").Value = '1.152'
$ws.Range("E35").Value = '  +4.55%  '
$ws.Range("D36").Value = '2.703'
$ws.Range("E36").Value = '  -1.02%  '
$ws.Range("D37").Value = '0.01880'
$ws.Range("E37").Value = '  +3.16%  '
$ws.Range("D38").Value = '2.644'
$ws.Range("E38").Value = '  -2.26%  '
$ws.Range("D39").Value = '0.9195'
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("D40").Value = '2.088'
$ws.Range("E40").Value = '  +0.88%  '
$ws.Range("D41").Value = '5.975'
$ws.Range("E41").Value = '  +1.04%  '
$ws.Range("D42").Value = '0.4356'
$ws.Range("E42").Value = '  -1.43%  '
$ws.Range("D43").Value = '106.16'
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("D45").Value = '7.643'
$ws.Range("E45").Value = '  +1.17%  '
$ws.Range("D46").Value = '1.604'
$ws.Range("E46").Value = '  +12.56%  '
$ws.Range("D47").Value = '0.1325'
$ws.Range("E47").Value = '  -3.83%  '
$ws.Range("D48").Value = '65.54'
$ws.Range("E48").Value = '  -13.19%  '
$ws.Range("D49").Value = '8.951'
$ws.Range("E49").Value = '  -1.22%  '
$ws.Range("D50").Value = '34.57'
$ws.Range("E50").Value = '  -2.74%  '
$ws.Range("D51").Value = '0.05714'
$ws.Range("E51").Value = '  -2.17%  '

# Restore default (General) style on column D so only the value/text type differs
$ws.Range("D2:D51").Style = "Normal"

